$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "config" sheet: append two new biofuel subtypes to the EBT_biofuels_list
#    (A15 = 15_solid_biomass_unallocated, A16 = 16_others_unallocated), using
#    the same black-font style ("s=6" in the original file) as the other
#    recently-added entries above them (A11:A14).
#    NB: write A16's new string first so the shared-string table gets the
#    same append order as the target file (64 = 16_others_unallocated,
#    65 = 15_solid_biomass_unallocated).
# ---------------------------------------------------------------------------
$wsConfig = $wb.Worksheets.Item("config")

$wsConfig.Range("A16").Value = "16_others_unallocated"
$wsConfig.Range("A16").Font.Color = 0

$wsConfig.Range("A15").Value = "15_solid_biomass_unallocated"
$wsConfig.Range("A15").Font.Color = 0

# ---------------------------------------------------------------------------
# 2) "simplified_economy_fuels" sheet: add a default method row, for every
#    economy, for each of the two new subfuels (mirrors the existing blocks
#    for the other subfuels already in the sheet).
# ---------------------------------------------------------------------------
$wsFuels = $wb.Worksheets.Item("simplified_economy_fuels")

$economies = @("00_MARS","01_AUS","02_BD","03_CDA","04_CHL","05_PRC","06_HKC","07_INA","08_JPN","09_ROK","10_MAS","11_MEX","12_NZ","13_PNG","14_PE","15_PHL","16_RUS","17_SGP","18_CT","19_THA","20_USA","21_VN")

$newFuels = @("16_others_unallocated","15_solid_biomass_unallocated")

$row = 310
foreach ($fuel in $newFuels) {
    foreach ($economy in $economies) {
        $wsFuels.Cells.Item($row, 1).Value = $economy
        $bCell = $wsFuels.Cells.Item($row, 2)
        $bCell.Value = $fuel
        $bCell.Font.Color = 0
        if ($economy -eq "00_MARS") {
            $wsFuels.Cells.Item($row, 3).Value = "satisfy_all_demand_with_domestic_production"
        } else {
            $wsFuels.Cells.Item($row, 3).Value = "satisfy_all_demand_with_domestic_production_EXACT"
        }
        $row = $row + 1
    }
}

# ---------------------------------------------------------------------------
# 3) View-state bookkeeping: the "config" sheet becomes the active tab /
#    selected cell A16, "simplified_economy_fuels" is left scrolled down with
#    B330 selected (and is no longer the active tab).
# ---------------------------------------------------------------------------
$wsFuels.Range("B330").Select()

$wsConfig.Activate()
$wsConfig.Range("A16").Select()
